# Updates cryptos list values (Price / Volume(1h)) to match the latest
# coinranking.com snapshot, and inserts a new "BabyDogeCoin" row at row 46,
# shifting Aptos..Algorand down by one row and dropping the former last
# row (Cronos), per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "0.9998", "1.000")
# are written with a leading apostrophe so Excel stores them as literal
# text (matching the inlineStr cells in the workbook) instead of coercing
# them to numeric values and dropping significant trailing/format digits.

# Row 2
$ws.Cells.Item(2, 4).Value = "29.404.81"
$ws.Cells.Item(2, 5).Value = "  -0.23%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.848.84"
$ws.Cells.Item(3, 5).Value = "  -0.05%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'0.9986"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'240.57"
$ws.Cells.Item(5, 5).Value = "  -0.52%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.6318"
$ws.Cells.Item(6, 5).Value = "  +0.50%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.9998"
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.07565"
$ws.Cells.Item(8, 5).Value = "  +0.70%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.2963"
$ws.Cells.Item(9, 5).Value = "  -0.33%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'24.48"
$ws.Cells.Item(10, 5).Value = "  +0.37%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.07706"
$ws.Cells.Item(11, 5).Value = "  -0.37%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.867.19"
$ws.Cells.Item(12, 5).Value = "  +0.68%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'5.004"
$ws.Cells.Item(13, 5).Value = "  +0.02%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.6854"
$ws.Cells.Item(14, 5).Value = "  -0.95%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.00001003"
$ws.Cells.Item(15, 5).Value = "  +3.10%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'83.04"
$ws.Cells.Item(16, 5).Value = "  -0.44%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.121.58"
$ws.Cells.Item(17, 5).Value = "  +0.17%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'6.165"
$ws.Cells.Item(18, 5).Value = "  -1.18%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "29.435.95"
$ws.Cells.Item(19, 5).Value = "  -0.25%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'228.64"
$ws.Cells.Item(20, 5).Value = "  -1.74%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'12.49"
$ws.Cells.Item(21, 5).Value = "  +0.06%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.02%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'7.572"
$ws.Cells.Item(23, 5).Value = "  -0.70%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'1.000"
$ws.Cells.Item(24, 5).Value = "  +0.00%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'158.41"
$ws.Cells.Item(25, 5).Value = "  +2.40%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.1398"
$ws.Cells.Item(26, 5).Value = "  +0.82%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'8.390"
$ws.Cells.Item(27, 5).Value = "  -0.75%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'17.68"
$ws.Cells.Item(28, 5).Value = "  -0.04%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.53%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +1.18%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.05696"
$ws.Cells.Item(31, 5).Value = "  -3.77%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'4.128"
$ws.Cells.Item(32, 5).Value = "  +0.60%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'4.028"
$ws.Cells.Item(33, 5).Value = "  +0.06%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.850"
$ws.Cells.Item(34, 5).Value = "  -2.24%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.98%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.7139"
$ws.Cells.Item(36, 5).Value = "  -0.86%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'2.590"
$ws.Cells.Item(37, 5).Value = "  +0.02%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "1.262.44"
$ws.Cells.Item(38, 5).Value = "  +2.07%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.01816"
$ws.Cells.Item(39, 5).Value = "  +1.30%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.61%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.9057"
$ws.Cells.Item(41, 5).Value = "  -0.07%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'6.171"
$ws.Cells.Item(42, 5).Value = "  +0.81%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'101.59"
$ws.Cells.Item(44, 5).Value = "  +0.17%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'66.37"
$ws.Cells.Item(45, 5).Value = "  -1.10%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(46, 4).Value = "'0.00000000120"
$ws.Cells.Item(46, 5).Value = "  +1.18%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Aptos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(47, 4).Value = "'7.079"
$ws.Cells.Item(47, 5).Value = "  -4.07%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "TheSandbox"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(48, 4).Value = "'0.4056"
$ws.Cells.Item(48, 5).Value = "  +0.49%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "'9.063"
$ws.Cells.Item(49, 5).Value = "  -0.73%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(50, 4).Value = "'1.684"
$ws.Cells.Item(50, 5).Value = "  -0.96%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).Value = "'0.1121"
$ws.Cells.Item(51, 5).Value = "  +0.43%  "
